# Revert unit test coverage: re-add row 39 (new sample record) to each of the
# four worksheets, restoring the dimension to A1:I39.

$wb = $excel.ActiveWorkbook

$dateSerial = 45825.46394675926
$dateFormat = "YYYY-MM-DD HH:MM:SS"

$rows = @{
    "MID_LFT_#1" = @{
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x7C"
        E = "0x07"
        F = 400
        G = 568631262647113000000000.0
        H = 380
        I = 7
    }
    "MID_LFT_#2" = @{
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x6C"
        E = "0x19"
        F = 380
        G = 568432987514711000000000.0
        H = 364
        I = 25
    }
    "MID_PLT_#1" = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x6B"
        E = "0x15"
        F = 110
        G = 568631262647113000000000.0
        H = 107
        I = 15
    }
    "MID_PLT_#2" = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x7F"
        E = "0x9"
        F = 130
        G = 568631262647113000000000.0
        H = 127
        I = 9
    }
}

foreach ($sheetName in $rows.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $rows[$sheetName]

    $ws.Range("A39").Value = $dateSerial
    $ws.Range("A39").NumberFormat = $dateFormat

    $ws.Range("B39").Value = $data.B
    $ws.Range("C39").Value = $data.C
    $ws.Range("D39").Value = $data.D
    $ws.Range("E39").Value = $data.E
    $ws.Range("F39").Value = $data.F
    $ws.Range("G39").Value = $data.G
    $ws.Range("H39").Value = $data.H
    $ws.Range("I39").Value = $data.I
}
